$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Wnt5a"
$ws.Cells.Item(2, 3).Value = "Lrp5"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.125615
$ws.Cells.Item(2, 8).Value = 0.25123
$ws.Cells.Item(2, 9).Value = 0.02647478672532295
$ws.Cells.Item(2, 10).Value = 0.01780700335556722
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 17.4294175
$ws.Cells.Item(2, 14).Value = 34.858835
$ws.Cells.Item(2, 15).Value = 0.4529581854295807
$ws.Cells.Item(2, 16).Value = 0.3776014560521451
$ws.Cells.Item(2, 17).Value = 2.1893962792625
$ws.Cells.Item(2, 18).Value = 8.75758511705
$ws.Cells.Item(2, 19).Value = 0.01199197135473744
$ws.Cells.Item(2, 20).Value = 0.006723950394987617

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Wnt5a"
$ws.Cells.Item(3, 3).Value = "Lrp5"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.125615
$ws.Cells.Item(3, 8).Value = 0.25123
$ws.Cells.Item(3, 9).Value = 0.02647478672532295
$ws.Cells.Item(3, 10).Value = 0.01780700335556722
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 3.914977
$ws.Cells.Item(3, 14).Value = 11.744931
$ws.Cells.Item(3, 15).Value = 0.1017429801035258
$ws.Cells.Item(3, 16).Value = 0.127224648983019
$ws.Cells.Item(3, 17).Value = 0.4917798358550001
$ws.Cells.Item(3, 18).Value = 2.95067901513
$ws.Cells.Item(3, 19).Value = 0.002693623699039621
$ws.Cells.Item(3, 20).Value = 0.002265489751351481

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Wnt5a"
$ws.Cells.Item(4, 3).Value = "Lrp5"
$ws.Cells.Item(4, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.125615
$ws.Cells.Item(4, 8).Value = 0.25123
$ws.Cells.Item(4, 9).Value = 0.02647478672532295
$ws.Cells.Item(4, 10).Value = 0.01780700335556722
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 5.397313
$ws.Cells.Item(4, 14).Value = 16.191939
$ws.Cells.Item(4, 15).Value = 0.1402661392829386
$ws.Cells.Item(4, 16).Value = 0.1753959862028526
$ws.Cells.Item(4, 17).Value = 0.6779834724950001
$ws.Cells.Item(4, 18).Value = 4.067900834970001
$ws.Cells.Item(4, 19).Value = 0.003713516122300243
$ws.Cells.Item(4, 20).Value = 0.003123276914867218

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Wnt5a"
$ws.Cells.Item(5, 3).Value = "Lrp5"
$ws.Cells.Item(5, 4).Value = "MuSCs"
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.125615
$ws.Cells.Item(5, 8).Value = 0.25123
$ws.Cells.Item(5, 9).Value = 0.02647478672532295
$ws.Cells.Item(5, 10).Value = 0.01780700335556722
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 5.69137
$ws.Cells.Item(5, 14).Value = 11.38274
$ws.Cells.Item(5, 15).Value = 0.1479081344978025
$ws.Cells.Item(5, 16).Value = 0.1233012863987851
$ws.Cells.Item(5, 17).Value = 0.71492144255
$ws.Cells.Item(5, 18).Value = 2.8596857702
$ws.Cells.Item(5, 19).Value = 0.003915836315769704
$ws.Cells.Item(5, 20).Value = 0.002195626420648922

# Row 6
$ws.Cells.Item(6, 1).Value = "ECs"
$ws.Cells.Item(6, 2).Value = "Wnt5a"
$ws.Cells.Item(6, 3).Value = "Lrp5"
$ws.Cells.Item(6, 4).Value = "Neutrophils"
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.125615
$ws.Cells.Item(6, 8).Value = 0.25123
$ws.Cells.Item(6, 9).Value = 0.02647478672532295
$ws.Cells.Item(6, 10).Value = 0.01780700335556722
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 1.988496
$ws.Cells.Item(6, 14).Value = 5.965488000000001
$ws.Cells.Item(6, 15).Value = 0.05167731737988258
$ws.Cells.Item(6, 16).Value = 0.06461997237892773
$ws.Cells.Item(6, 17).Value = 0.24978492504
$ws.Cells.Item(6, 18).Value = 1.49870955024
$ws.Cells.Item(6, 19).Value = 0.001368145956169217
$ws.Cells.Item(6, 20).Value = 0.001150688064988227

# Row 7
$ws.Cells.Item(7, 1).Value = "ECs"
$ws.Cells.Item(7, 2).Value = "Wnt5a"
$ws.Cells.Item(7, 3).Value = "Lrp5"
$ws.Cells.Item(7, 4).Value = "Resolving-Mac"
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.125615
$ws.Cells.Item(7, 8).Value = 0.25123
$ws.Cells.Item(7, 9).Value = 0.02647478672532295
$ws.Cells.Item(7, 10).Value = 0.01780700335556722
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 4.057513666666667
$ws.Cells.Item(7, 14).Value = 12.172541
$ws.Cells.Item(7, 15).Value = 0.1054472433062699
$ws.Cells.Item(7, 16).Value = 0.1318566499842704
$ws.Cells.Item(7, 17).Value = 0.5096845792383334
$ws.Cells.Item(7, 18).Value = 3.05810747543
$ws.Cells.Item(7, 19).Value = 0.002791693277306734
$ws.Cells.Item(7, 20).Value = 0.002347971808723756

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Wnt5a"
$ws.Cells.Item(8, 3).Value = "Lrp5"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 4.619088000000001
$ws.Cells.Item(8, 8).Value = 13.857264
$ws.Cells.Item(8, 9).Value = 0.9735252132746771
$ws.Cells.Item(8, 10).Value = 0.9821929966444328
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 17.4294175
$ws.Cells.Item(8, 14).Value = 34.858835
$ws.Cells.Item(8, 15).Value = 0.4529581854295807
$ws.Cells.Item(8, 16).Value = 0.3776014560521451
$ws.Cells.Item(8, 17).Value = 80.50801322124
$ws.Cells.Item(8, 18).Value = 483.04807932744
$ws.Cells.Item(8, 19).Value = 0.4409662140748433
$ws.Cells.Item(8, 20).Value = 0.3708775056571575

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Wnt5a"
$ws.Cells.Item(9, 3).Value = "Lrp5"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 4.619088000000001
$ws.Cells.Item(9, 8).Value = 13.857264
$ws.Cells.Item(9, 9).Value = 0.9735252132746771
$ws.Cells.Item(9, 10).Value = 0.9821929966444328
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 3.914977
$ws.Cells.Item(9, 14).Value = 11.744931
$ws.Cells.Item(9, 15).Value = 0.1017429801035258
$ws.Cells.Item(9, 16).Value = 0.127224648983019
$ws.Cells.Item(9, 17).Value = 18.083623280976
$ws.Cells.Item(9, 18).Value = 162.752609528784
$ws.Cells.Item(9, 19).Value = 0.09904935640448616
$ws.Cells.Item(9, 20).Value = 0.1249591592316675

# Row 10
$ws.Cells.Item(10, 1).Value = "FAPs"
$ws.Cells.Item(10, 2).Value = "Wnt5a"
$ws.Cells.Item(10, 3).Value = "Lrp5"
$ws.Cells.Item(10, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 4.619088000000001
$ws.Cells.Item(10, 8).Value = 13.857264
$ws.Cells.Item(10, 9).Value = 0.9735252132746771
$ws.Cells.Item(10, 10).Value = 0.9821929966444328
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 5.397313
$ws.Cells.Item(10, 14).Value = 16.191939
$ws.Cells.Item(10, 15).Value = 0.1402661392829386
$ws.Cells.Item(10, 16).Value = 0.1753959862028526
$ws.Cells.Item(10, 17).Value = 24.93066371054401
$ws.Cells.Item(10, 18).Value = 224.375973394896
$ws.Cells.Item(10, 19).Value = 0.1365526231606383
$ws.Cells.Item(10, 20).Value = 0.1722727092879854

# Row 11
$ws.Cells.Item(11, 1).Value = "FAPs"
$ws.Cells.Item(11, 2).Value = "Wnt5a"
$ws.Cells.Item(11, 3).Value = "Lrp5"
$ws.Cells.Item(11, 4).Value = "MuSCs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 4.619088000000001
$ws.Cells.Item(11, 8).Value = 13.857264
$ws.Cells.Item(11, 9).Value = 0.9735252132746771
$ws.Cells.Item(11, 10).Value = 0.9821929966444328
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 5.69137
$ws.Cells.Item(11, 14).Value = 11.38274
$ws.Cells.Item(11, 15).Value = 0.1479081344978025
$ws.Cells.Item(11, 16).Value = 0.1233012863987851
$ws.Cells.Item(11, 17).Value = 26.28893887056
$ws.Cells.Item(11, 18).Value = 157.73363322336
$ws.Cells.Item(11, 19).Value = 0.1439922981820328
$ws.Cells.Item(11, 20).Value = 0.1211056599781362

# Row 12
$ws.Cells.Item(12, 1).Value = "FAPs"
$ws.Cells.Item(12, 2).Value = "Wnt5a"
$ws.Cells.Item(12, 3).Value = "Lrp5"
$ws.Cells.Item(12, 4).Value = "Neutrophils"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 4.619088000000001
$ws.Cells.Item(12, 8).Value = 13.857264
$ws.Cells.Item(12, 9).Value = 0.9735252132746771
$ws.Cells.Item(12, 10).Value = 0.9821929966444328
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 1.988496
$ws.Cells.Item(12, 14).Value = 5.965488000000001
$ws.Cells.Item(12, 15).Value = 0.05167731737988258
$ws.Cells.Item(12, 16).Value = 0.06461997237892773
$ws.Cells.Item(12, 17).Value = 9.185038011648002
$ws.Cells.Item(12, 18).Value = 82.665342104832
$ws.Cells.Item(12, 19).Value = 0.05030917142371337
$ws.Cells.Item(12, 20).Value = 0.06346928431393951

# Row 13
$ws.Cells.Item(13, 1).Value = "FAPs"
$ws.Cells.Item(13, 2).Value = "Wnt5a"
$ws.Cells.Item(13, 3).Value = "Lrp5"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 4.619088000000001
$ws.Cells.Item(13, 8).Value = 13.857264
$ws.Cells.Item(13, 9).Value = 0.9735252132746771
$ws.Cells.Item(13, 10).Value = 0.9821929966444328
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 4.057513666666667
$ws.Cells.Item(13, 14).Value = 12.172541
$ws.Cells.Item(13, 15).Value = 0.1054472433062699
$ws.Cells.Item(13, 16).Value = 0.1318566499842704
$ws.Cells.Item(13, 17).Value = 18.742012687536
$ws.Cells.Item(13, 18).Value = 168.678114187824
$ws.Cells.Item(13, 19).Value = 0.1026555500289631
$ws.Cells.Item(13, 20).Value = 0.1295086781755467
